# Re-sort the calibration data rows (A2:D12) in ascending order of time (column A).
# This reflects the "Performed calibration of the needle" commit, where the
# curvature samples end up chronologically ordered by the time (s) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortedData = @(
    @(57034.363191, -0.0000032422421282, -0.0000087073497606, -0.000024915333007),
    @(57046.031192, -0.000066157834425, -0.000096710380647, -0.00010499781355),
    @(57059.299192, -0.000073728495201, -0.00022851370616, -0.00015041768697),
    @(57070.563193, -0.000081532374738, -0.00036044991622, -0.0001933029678),
    @(57081.963193, -0.0001228416, -0.0004886953, -0.0002440999),
    @(57091.831194, -0.0001927213, -0.0006102247, -0.0003006622),
    @(57102.299194, -0.000146344, -0.0004822301, -0.0002532177),
    @(57112.163195, -0.000096978917629, -0.00035120831059, -0.00020319603157),
    @(57123.231195, -0.000044424733057, -0.00022768781771, -0.00014148830698),
    @(57136.031196, -0.000026573538511, -0.000073390219873, -0.00010844868675),
    @(57147.963197, -0.0000092459399218, -0.000013298685729, -0.000025306092296)
)

$row = 2
foreach ($r in $sortedData) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$wb.Save()
